$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H39").Value = 205.625
$ws.Range("I39").Value = 82.42856999999999
$ws.Range("J39").Value = 301.44446
$ws.Range("K39").Value = 247.28571
$ws.Range("L39").Value = 904.33338
$ws.Range("M39").Value = 48.71429000000001
$ws.Range("N39").Value = -1496.33338
$ws.Range("H129").Value = 1053.6207
$ws.Range("I129").Value = 295.66666
$ws.Range("J129").Value = 1094.9636
$ws.Range("K129").Value = 886.9999799999999
$ws.Range("L129").Value = 3284.8908
$ws.Range("M129").Value = 4113.00002
$ws.Range("N129").Value = -13284.8908
$ws.Range("H133").Value = 43407.5
$ws.Range("J133").Value = 43407.5
$ws.Range("L133").Value = 43407.5
$ws.Range("N133").Value = -53527.5
$ws.Range("H138").Value = 10163535
$ws.Range("I138").Value = 3099218.2
$ws.Range("J138").Value = 13891924
$ws.Range("K138").Value = 9297654.600000001
$ws.Range("L138").Value = 41675772
$ws.Range("M138").Value = -9292514.600000001
$ws.Range("N138").Value = -41686052
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 22275.117
$ws.Range("I32").Value = 4790.857
$ws.Range("K32").Value = 4790.857
$ws.Range("M32").Value = -4503.857
$ws.Range("H132").Value = 3264.1143
$ws.Range("I132").Value = 2722.76
$ws.Range("J132").Value = 4617.5
$ws.Range("K132").Value = 8168.280000000001
$ws.Range("L132").Value = 13852.5
$ws.Range("M132").Value = -5638.280000000001
$ws.Range("N132").Value = -18912.5
$ws.Range("H139").Value = 51833.875
$ws.Range("J139").Value = 51833.875
$ws.Range("L139").Value = 51833.875
$ws.Range("N139").Value = -62113.875
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H129").Value = 44499.5
$ws.Range("J129").Value = 44499.5
$ws.Range("L129").Value = 44499.5
$ws.Range("N129").Value = -54499.5
$ws.Range("H133").Value = 39780
$ws.Range("J133").Value = 39780
$ws.Range("L133").Value = 39780
$ws.Range("N133").Value = -49900
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 2658.0588
$ws.Range("I99").Value = 1998.6666
$ws.Range("J99").Value = 3399.875
$ws.Range("K99").Value = 1998.6666
$ws.Range("L99").Value = 3399.875
$ws.Range("M99").Value = -500.6666
$ws.Range("N99").Value = -6395.875
$ws.Range("H126").Value = 2658.0588
$ws.Range("I126").Value = 1998.6666
$ws.Range("J126").Value = 3399.875
$ws.Range("K126").Value = 5995.9998
$ws.Range("L126").Value = 10199.625
$ws.Range("M126").Value = -3525.9998
$ws.Range("N126").Value = -15139.625
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1137.2632
$ws.Range("I5").Value = 693.53845
$ws.Range("J5").Value = 1509.4193
$ws.Range("K5").Value = 2080.61535
$ws.Range("L5").Value = 4528.257900000001
$ws.Range("M5").Value = -1968.61535
$ws.Range("N5").Value = -4752.257900000001
$ws.Range("H50").Value = 447.84616
$ws.Range("I50").Value = 370.33334
$ws.Range("J50").Value = 514.2857
$ws.Range("K50").Value = 1111.00002
$ws.Range("L50").Value = 1542.8571
$ws.Range("M50").Value = -630.0000199999999
$ws.Range("N50").Value = -2504.8571
$ws.Range("H53").Value = 447.84616
$ws.Range("I53").Value = 370.33334
$ws.Range("J53").Value = 514.2857
$ws.Range("K53").Value = 1111.00002
$ws.Range("L53").Value = 1542.8571
$ws.Range("M53").Value = -630.0000199999999
$ws.Range("N53").Value = -2504.8571
$ws.Range("H107").Value = 384.1111
$ws.Range("J107").Value = 341.9524
$ws.Range("L107").Value = 1025.8572
$ws.Range("N107").Value = -4865.8572
$ws.Range("H119").Value = 1144.625
$ws.Range("I119").Value = 592.8333
$ws.Range("K119").Value = 1778.4999
$ws.Range("M119").Value = 3059.5001
$ws.Range("H131").Value = 6804237
$ws.Range("I131").Value = 608
$ws.Range("J131").Value = 7577376.5
$ws.Range("K131").Value = 1824
$ws.Range("L131").Value = 22732129.5
$ws.Range("M131").Value = 3216
$ws.Range("N131").Value = -22742209.5
$ws.Range("H133").Value = 8928.888999999999
$ws.Range("H135").Value = 1137.2632
$ws.Range("I135").Value = 693.53845
$ws.Range("J135").Value = 1509.4193
$ws.Range("K135").Value = 6241.84605
$ws.Range("L135").Value = 13584.7737
$ws.Range("M135").Value = -3706.84605
$ws.Range("N135").Value = -18654.7737
$ws.Range("H140").Value = 7840.294
$ws.Range("I140").Value = 11278.5
$ws.Range("J140").Value = 2928.5715
$ws.Range("K140").Value = 33835.5
$ws.Range("L140").Value = 8785.7145
$ws.Range("M140").Value = -28655.5
$ws.Range("N140").Value = -19145.7145
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4396.2666
$ws.Range("I70").Value = 4343.923
$ws.Range("J70").Value = 4736.5
$ws.Range("K70").Value = 4343.923
$ws.Range("L70").Value = 4736.5
$ws.Range("M70").Value = -4073.923
$ws.Range("N70").Value = -5276.5
$ws.Range("H73").Value = 4396.2666
$ws.Range("I73").Value = 4343.923
$ws.Range("J73").Value = 4736.5
$ws.Range("K73").Value = 4343.923
$ws.Range("L73").Value = 4736.5
$ws.Range("M73").Value = -3407.923
$ws.Range("N73").Value = -6608.5
$ws.Range("H122").Value = 2925
$ws.Range("I122").Value = 3000
$ws.Range("J122").Value = 2900
$ws.Range("K122").Value = 9000
$ws.Range("L122").Value = 8700
$ws.Range("M122").Value = -6550
$ws.Range("N122").Value = -13600
$ws.Range("H138").Value = 64861.285
$ws.Range("J138").Value = 64861.285
$ws.Range("L138").Value = 64861.285
$ws.Range("N138").Value = -75141.285
$ws.Range("H139").Value = 41269
$ws.Range("J139").Value = 41269
$ws.Range("L139").Value = 41269
$ws.Range("N139").Value = -51549
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3140
$ws.Range("I7").Value = 1866.6666
$ws.Range("J7").Value = 3458.3333
$ws.Range("K7").Value = 1866.6666
$ws.Range("L7").Value = 3458.3333
$ws.Range("M7").Value = -1754.6666
$ws.Range("N7").Value = -3682.3333
$ws.Range("H40").Value = 2503.8845
$ws.Range("I40").Value = 1707.8462
$ws.Range("J40").Value = 3299.923
$ws.Range("K40").Value = 1707.8462
$ws.Range("L40").Value = 3299.923
$ws.Range("M40").Value = -1571.8462
$ws.Range("N40").Value = -3571.923
$ws.Range("H46").Value = 1461.5385
$ws.Range("I46").Value = 1066.6666
$ws.Range("J46").Value = 1580
$ws.Range("K46").Value = 1066.6666
$ws.Range("L46").Value = 1580
$ws.Range("M46").Value = -878.6666
$ws.Range("N46").Value = -1956
$ws.Range("H126").Value = 3140
$ws.Range("I126").Value = 1866.6666
$ws.Range("J126").Value = 3458.3333
$ws.Range("K126").Value = 5599.9998
$ws.Range("L126").Value = 10374.9999
$ws.Range("M126").Value = -3129.9998
$ws.Range("N126").Value = -15314.9999
$ws.Range("H132").Value = 4568.091
$ws.Range("I132").Value = 2741.1428
$ws.Range("J132").Value = 7765.25
$ws.Range("K132").Value = 8223.428400000001
$ws.Range("L132").Value = 23295.75
$ws.Range("M132").Value = -5693.428400000001
$ws.Range("N132").Value = -28355.75
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3500.838
$ws.Range("I132").Value = 3488.2666
$ws.Range("J132").Value = 3554.7144
$ws.Range("K132").Value = 10464.7998
$ws.Range("L132").Value = 10664.1432
$ws.Range("M132").Value = -7934.799800000001
$ws.Range("N132").Value = -15724.1432
